$d = $word.ActiveDocument

# The only visible content change in the diff: the "cd" path in the code block
# changes from ".../practicals/prac5/examples/" to ".../practicals/prac5/code"
$d.Content.Find.Execute(
    "CWM-in-HPC-and-Scientific-Computing-2020/practicals/prac5/examples/",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "CWM-in-HPC-and-Scientific-Computing-2020/practicals/prac5/code",
    2)
